$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 22 numeric updates
$ws1.Range("C22").Value = 1555.2
$ws1.Range("D22").Value = 798.3099999999999
$ws1.Range("L22").Value = 380.16

# Row 55 "x de 53" summary labels
$ws1.Range("C55").Value = "3 de 53"
$ws1.Range("D55").Value = "9 de 53"
$ws1.Range("L55").Value = "5 de 53"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F22").Value = 2733.67
$ws2.Range("F55").Value = 54694.17

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D2").Value = 5365.44
$ws3.Range("E2").Value = 454.5600000000004
$ws3.Range("F2").Value = 0.9218969072164948

$ws3.Range("D3").Value = 26253.99
$ws3.Range("E3").Value = -12525.99
$ws3.Range("F3").Value = 1.912440996503497

$ws3.Range("D15").Value = 10066.44
$ws3.Range("E15").Value = 5623.559999999999
$ws3.Range("F15").Value = 0.6415831739961759

$ws3.Range("D19").Value = 56087.77
$ws3.Range("E19").Value = 34875.55899999999
$ws3.Range("F19").Value = 0.6165975961587774
